$wb = $excel.ActiveWorkbook

# --- Rename sheets (task order id timestamps updated) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16504777862383795"
$wb.Worksheets.Item(2).Name = "NB_TO-1650477787817412"
$wb.Worksheets.Item(3).Name = "RS_TO-16504777878183773"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504777878663797"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650477787929412"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504777862083795.csv"
$ws1.Range("B3").Value = "GNG_stims-16504777862213771.csv"
$ws1.Range("B4").Value = "go_stims-1650477786222379.csv"
$ws1.Range("B5").Value = "GNG_stims-16504777862373767.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16504777872603772.csv"
$ws2.Range("B3").Value = "ZB-match_0-16504777863464117.csv"
$ws2.Range("B4").Value = "TB-16504777877974117.csv"
$ws2.Range("B5").Value = "ZB-match_8-16504777864104111.csv"
$ws2.Range("B6").Value = "OB-165047778647541.csv"
$ws2.Range("B7").Value = "ZB-match_3-16504777862874067.csv"
$ws2.Range("B8").Value = "TB-16504777876743777.csv"
$ws2.Range("B9").Value = "OB-16504777868564103.csv"
$ws2.Range("B10").Value = "OB-16504777865083745.csv"

# --- Sheet 3 (RS) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504777878334122.csv"
$ws4.Range("B3").Value = "ZM_stims-16504777878193784.csv"
$ws4.Range("B4").Value = "MM_stims-16504777878494089.csv"
$ws4.Range("B5").Value = "ZM_stims-16504777878334122.csv"
$ws4.Range("B6").Value = "MM_stims-1650477787865376.csv"
$ws4.Range("B7").Value = "ZM_stims-16504777878494089.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1650477787914377.csv"
$ws5.Range("B3").Value = "SAT_stims-16504777878683748.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504777878974116.csv"
$ws5.Range("B5").Value = "SAT_stims-16504777878814113.csv"
